$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$data = @(
    @(4,2,6,0),
    @(3,2,3,1),
    @(4,1,5,2),
    @(3,2,3,1),
    @(5,0,6,2),
    @(4,2,3,1),
    @(4,0,4,2),
    @(6,2,6,0),
    @(6,2,6,0),
    @(3,2,3,1),
    @(3,1,4,2),
    @(5,0,5,2),
    @(5,3,2,0),
    @(5,2,6,0),
    @(4,1,2,2),
    @(3,1,4,2),
    @(5,2,6,0),
    @(3,0,4,3),
    @(3,3,4,0),
    @(7,0,5,3),
    @(6,2,6,1),
    @(5,2,6,1),
    @(2,2,4,1),
    @(4,1,3,2),
    @(3,2,3,1)
)

$startRow = 1313
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
}

$nextRow = $startRow + $data.Count
$ws.Range("A$nextRow").Select() | Out-Null
